$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Step 1: update the date serial in A1 (was 45308 -> 45309, i.e. +1 day)
$ws.Range("A1").Value = 45309

# Step 2: update price values in the "altas" (TUA-1xx) price table
$ws.Range("D23").Value = 11284.427
$ws.Range("D24").Value = 15555.424
$ws.Range("D25").Value = 22119.28
$ws.Range("D26").Value = 32999.091
$ws.Range("D27").Value = 59524.242
$ws.Range("D28").Value = 77057.803

# Update price values in the "bajas" (TU-1xx) price table
$ws.Range("D36").Value = 6250
$ws.Range("D37").Value = 9150
